$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = 3
$ws.Range("D2").Value = "pqr"
$ws.Range("D3").Value = "pune"

$ws.Range("D4").Select() | Out-Null
